$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Cells.Item(2, 4).Value = 7989
$ws.Cells.Item(2, 5).Value = 234
$ws.Cells.Item(2, 6).Value = 234
$ws.Cells.Item(2, 7).Value = 50
$ws.Cells.Item(2, 8).Value = 11
$ws.Cells.Item(2, 9).Value = 6
$ws.Cells.Item(2, 10).Value = 5
$ws.Cells.Item(2, 11).Value = 6019
$ws.Cells.Item(2, 12).Value = 4340
$ws.Cells.Item(2, 13).Value = 1679
$ws.Cells.Item(2, 14).Value = 1659
$ws.Cells.Item(2, 15).Value = 19
$ws.Cells.Item(2, 16).Value = 196
$ws.Cells.Item(2, 17).Value = 119
$ws.Cells.Item(2, 18).Value = -272
$ws.Cells.Item(2, 19).Value = 224
$ws.Cells.Item(2, 20).Value = 228
$ws.Cells.Item(2, 21).Value = -109
$ws.Cells.Item(2, 22).Value = 3102
$ws.Cells.Item(2, 23).Value = 2.93
$ws.Cells.Item(2, 24).Value = 0.14
$ws.Cells.Item(2, 25).Value = 0.34
$ws.Cells.Item(2, 26).Value = 0.19
$ws.Cells.Item(2, 27).Value = 258.58
$ws.Cells.Item(2, 28).Value = 736.29
$ws.Cells.Item(2, 29).Value = 14
$ws.Cells.Item(2, 30).Value = 402.19
$ws.Cells.Item(2, 31).Value = 4491
$ws.Cells.Item(2, 32).Value = 1.3
$ws.Cells.Item(2, 33).Value = 25
$ws.Cells.Item(2, 34).Value = 0.43
$ws.Cells.Item(2, 35).Value = 166.3
$ws.Cells.Item(2, 36).Value = 39258141
$ws.Cells.Item(3, 4).Value = 8426
$ws.Cells.Item(3, 5).Value = 183
$ws.Cells.Item(3, 6).Value = 183
$ws.Cells.Item(3, 7).Value = -138
$ws.Cells.Item(3, 8).Value = -173
$ws.Cells.Item(3, 9).Value = -176
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(3, 11).Value = 6112
$ws.Cells.Item(3, 12).Value = 4554
$ws.Cells.Item(3, 13).Value = 1559
$ws.Cells.Item(3, 14).Value = 1537
$ws.Cells.Item(3, 15).Value = 22
$ws.Cells.Item(3, 16).Value = 199
$ws.Cells.Item(3, 17).Value = 12
$ws.Cells.Item(3, 18).Value = -512
$ws.Cells.Item(3, 19).Value = 492
$ws.Cells.Item(3, 20).Value = 333
$ws.Cells.Item(3, 21).Value = -320
$ws.Cells.Item(3, 22).Value = 3493
$ws.Cells.Item(3, 23).Value = 2.17
$ws.Cells.Item(3, 24).Value = -2.06
$ws.Cells.Item(3, 25).Value = -11
$ws.Cells.Item(3, 26).Value = -2.86
$ws.Cells.Item(3, 27).Value = 292.11
$ws.Cells.Item(3, 28).Value = 634.03
$ws.Cells.Item(3, 29).Value = -445
$ws.Cells.Item(3, 30).Value = -9.09
$ws.Cells.Item(3, 31).Value = 4111
$ws.Cells.Item(3, 32).Value = 0.98
$ws.Cells.Item(3, 33).Value = 20
$ws.Cells.Item(3, 34).Value = 0.49
$ws.Cells.Item(3, 35).Value = -4.25
$ws.Cells.Item(3, 36).Value = 39706670
$ws.Cells.Item(4, 4).Value = 8647
$ws.Cells.Item(4, 5).Value = -486
$ws.Cells.Item(4, 6).Value = -486
$ws.Cells.Item(4, 7).Value = -637
$ws.Cells.Item(4, 8).Value = -533
$ws.Cells.Item(4, 9).Value = -537
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(4, 11).Value = 5537
$ws.Cells.Item(4, 12).Value = 4405
$ws.Cells.Item(4, 13).Value = 1133
$ws.Cells.Item(4, 14).Value = 1107
$ws.Cells.Item(4, 15).Value = 26
$ws.Cells.Item(4, 16).Value = 216
$ws.Cells.Item(4, 17).Value = 389
$ws.Cells.Item(4, 18).Value = -134
$ws.Cells.Item(4, 19).Value = -61
$ws.Cells.Item(4, 20).Value = 183
$ws.Cells.Item(4, 21).Value = 206
$ws.Cells.Item(4, 22).Value = 3402
$ws.Cells.Item(4, 23).Value = -5.62
$ws.Cells.Item(4, 24).Value = -6.17
$ws.Cells.Item(4, 25).Value = -40.64
$ws.Cells.Item(4, 26).Value = -9.15
$ws.Cells.Item(4, 27).Value = 388.92
$ws.Cells.Item(4, 28).Value = 373.63
$ws.Cells.Item(4, 29).Value = -1342
$ws.Cells.Item(4, 30).Value = -2.9
$ws.Cells.Item(4, 31).Value = 2703
$ws.Cells.Item(4, 32).Value = 1.44
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 43269712
$ws.Cells.Item(5, 4).Value = 9222
$ws.Cells.Item(5, 5).Value = 237
$ws.Cells.Item(5, 6).Value = 237
$ws.Cells.Item(5, 7).Value = 152
$ws.Cells.Item(5, 8).Value = 99
$ws.Cells.Item(5, 9).Value = 96
$ws.Cells.Item(5, 10).Value = 3
$ws.Cells.Item(5, 11).Value = 5312
$ws.Cells.Item(5, 12).Value = 3942
$ws.Cells.Item(5, 13).Value = 1370
$ws.Cells.Item(5, 14).Value = 1345
$ws.Cells.Item(5, 15).Value = 24
$ws.Cells.Item(5, 16).Value = 246
$ws.Cells.Item(5, 17).Value = 213
$ws.Cells.Item(5, 18).Value = -156
$ws.Cells.Item(5, 19).Value = -148
$ws.Cells.Item(5, 20).Value = 85
$ws.Cells.Item(5, 21).Value = 129
$ws.Cells.Item(5, 22).Value = 2900
$ws.Cells.Item(5, 23).Value = 2.57
$ws.Cells.Item(5, 24).Value = 1.07
$ws.Cells.Item(5, 25).Value = 7.86
$ws.Cells.Item(5, 26).Value = 1.83
$ws.Cells.Item(5, 27).Value = 287.87
$ws.Cells.Item(5, 28).Value = 431.22
$ws.Cells.Item(5, 29).Value = 209
$ws.Cells.Item(5, 30).Value = 16.8
$ws.Cells.Item(5, 31).Value = 2864
$ws.Cells.Item(5, 32).Value = 1.23
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 49283000
$ws.Cells.Item(6, 4).Value = 9721
$ws.Cells.Item(6, 5).Value = 325
$ws.Cells.Item(6, 6).Value = 325
$ws.Cells.Item(6, 7).Value = 163
$ws.Cells.Item(6, 8).Value = 104
$ws.Cells.Item(6, 9).Value = 101
$ws.Cells.Item(6, 11).Value = 5565
$ws.Cells.Item(6, 12).Value = 4055
$ws.Cells.Item(6, 13).Value = 1510
$ws.Cells.Item(6, 14).Value = 1483
$ws.Cells.Item(6, 16).Value = 250
$ws.Cells.Item(6, 17).Value = 146
$ws.Cells.Item(6, 18).Value = -114
$ws.Cells.Item(6, 19).Value = 25
$ws.Cells.Item(6, 20).Value = 137
$ws.Cells.Item(6, 21).Value = 10
$ws.Cells.Item(6, 22).Value = 2912
$ws.Cells.Item(6, 23).Value = 3.35
$ws.Cells.Item(6, 24).Value = 1.07
$ws.Cells.Item(6, 25).Value = 7.14
$ws.Cells.Item(6, 26).Value = 1.91
$ws.Cells.Item(6, 27).Value = 268.57
$ws.Cells.Item(6, 28).Value = 463.54
$ws.Cells.Item(6, 29).Value = 202
$ws.Cells.Item(6, 30).Value = 12.67
$ws.Cells.Item(6, 31).Value = 3161
$ws.Cells.Item(6, 32).Value = 0.8100000000000001
$ws.Cells.Item(6, 33).Value = 20
$ws.Cells.Item(6, 34).Value = 0.78
$ws.Cells.Item(6, 35).Value = 9.300000000000001
$ws.Cells.Item(6, 36).Value = 49911742
$ws.Cells.Item(7, 4).Value = 9923
$ws.Cells.Item(7, 5).Value = 447
$ws.Cells.Item(7, 7).Value = 242
$ws.Cells.Item(7, 8).Value = 199
$ws.Cells.Item(7, 9).Value = 195
$ws.Cells.Item(7, 11).Value = 6540
$ws.Cells.Item(7, 12).Value = 4882
$ws.Cells.Item(7, 13).Value = 1656
$ws.Cells.Item(7, 14).Value = 1630
$ws.Cells.Item(7, 16).Value = 250
$ws.Cells.Item(7, 17).Value = 356
$ws.Cells.Item(7, 18).Value = -198
$ws.Cells.Item(7, 19).Value = 756
$ws.Cells.Item(7, 20).Value = 192
$ws.Cells.Item(7, 21).Value = 153
$ws.Cells.Item(7, 23).Value = 4.5
$ws.Cells.Item(7, 24).Value = 2
$ws.Cells.Item(7, 25).Value = 12.53
$ws.Cells.Item(7, 26).Value = 3.29
$ws.Cells.Item(7, 27).Value = 294.75
$ws.Cells.Item(7, 29).Value = 391
$ws.Cells.Item(7, 30).Value = 6.02
$ws.Cells.Item(7, 31).Value = 3503
$ws.Cells.Item(7, 32).Value = 0.67
$ws.Cells.Item(7, 33).Value = 20
$ws.Cells.Item(7, 34).Value = 0.85
$ws.Cells.Item(7, 35).Value = 5.12
$ws.Cells.Item(8, 4).Value = 10253
$ws.Cells.Item(8, 5).Value = 536
$ws.Cells.Item(8, 7).Value = 414
$ws.Cells.Item(8, 8).Value = 318
$ws.Cells.Item(8, 9).Value = 313
$ws.Cells.Item(8, 11).Value = 6930
$ws.Cells.Item(8, 12).Value = 4963
$ws.Cells.Item(8, 13).Value = 1967
$ws.Cells.Item(8, 14).Value = 1938
$ws.Cells.Item(8, 16).Value = 250
$ws.Cells.Item(8, 17).Value = 471
$ws.Cells.Item(8, 18).Value = -266
$ws.Cells.Item(8, 19).Value = 76
$ws.Cells.Item(8, 20).Value = 240
$ws.Cells.Item(8, 21).Value = 67
$ws.Cells.Item(8, 23).Value = 5.23
$ws.Cells.Item(8, 24).Value = 3.1
$ws.Cells.Item(8, 25).Value = 17.56
$ws.Cells.Item(8, 26).Value = 4.73
$ws.Cells.Item(8, 27).Value = 252.31
$ws.Cells.Item(8, 29).Value = 628
$ws.Cells.Item(8, 30).Value = 3.74
$ws.Cells.Item(8, 31).Value = 4166
$ws.Cells.Item(8, 32).Value = 0.5600000000000001
$ws.Cells.Item(8, 33).Value = 20
$ws.Cells.Item(8, 34).Value = 0.85
$ws.Cells.Item(8, 35).Value = 3.19
$ws.Cells.Item(9, 4).Value = 11021
$ws.Cells.Item(9, 5).Value = 575
$ws.Cells.Item(9, 7).Value = 482
$ws.Cells.Item(9, 8).Value = 369
$ws.Cells.Item(9, 9).Value = 367
$ws.Cells.Item(9, 11).Value = 7530
$ws.Cells.Item(9, 12).Value = 5161
$ws.Cells.Item(9, 13).Value = 2369
$ws.Cells.Item(9, 14).Value = 2343
$ws.Cells.Item(9, 16).Value = 250
$ws.Cells.Item(9, 17).Value = 484
$ws.Cells.Item(9, 18).Value = -190
$ws.Cells.Item(9, 19).Value = 146
$ws.Cells.Item(9, 20).Value = 160
$ws.Cells.Item(9, 21).Value = 209
$ws.Cells.Item(9, 23).Value = 5.22
$ws.Cells.Item(9, 24).Value = 3.35
$ws.Cells.Item(9, 25).Value = 17.16
$ws.Cells.Item(9, 26).Value = 5.1
$ws.Cells.Item(9, 27).Value = 217.86
$ws.Cells.Item(9, 29).Value = 736
$ws.Cells.Item(9, 30).Value = 3.19
$ws.Cells.Item(9, 31).Value = 5035
$ws.Cells.Item(9, 32).Value = 0.47
$ws.Cells.Item(9, 33).Value = 20
$ws.Cells.Item(9, 34).Value = 2.72